$wb = $excel.ActiveWorkbook

# The edit targets the "Repayment Schedule" sheet (already the active
# sheet / selected tab in this workbook).
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Activate()

# Insert a new blank column before column N. This shifts the existing
# "Late" column (N -> O) and "Outstanding" column (P -> Q) one slot to
# the right, making room for the new "Variable Instalments" column.
$ws.Columns("N:N").Insert() | Out-Null

# Matches the author's final cursor position after the edit.
$ws.Range("T5").Select() | Out-Null
